$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.138157486915588
$ws.Range("B1").Value = 2.227763891220093
$ws.Range("C1").Value = 10.84977531433105
$ws.Range("D1").Value = 2.31383228302002
$ws.Range("E1").Value = 1.277617573738098
